$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.65"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.320"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.500"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.129"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8181"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8669"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01008"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1375"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07020"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03214"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02897"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09402"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.755"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001520"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04712"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006176"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001236"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003849"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008797"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.535"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.139"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3174"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1330"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1328"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003013"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03710"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006398"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1055"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002216"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008650"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005258"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3883"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002253"
